# Pathfinder dark knight class.docx - wording changes (content unchanged in effect)
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "awakened" -> "awaken" in the class flavor-text intro paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "straddle the line of undeath and awakened many unusual",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "straddle the line of undeath and awaken many unusual", 2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) First stamina-drain ability (the 1d4 stamina / "fatigued" one): insert
#    " as the duration" right before the period that ends the first
#    sentence (right after "...recovers <fraction> as much stamina").
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute(
    "as much stamina. Successful save causes them to loose half as much stamina or be fatigued",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$insPos1 = $rng1.Start + "as much stamina".Length
$ip1 = $d.Range($insPos1, $insPos1)
$ip1.InsertAfter(" as the duration")

# ---------------------------------------------------------------------------
# 3) Second stamina-drain ability (the 1d6 stamina / "exhausted" one):
#    insert " as the duration" the same way, right before the period.
#    The document's lone "_GoBack" bookmark (currently sitting in the
#    "Vicious" paragraph, between the "T" and "his bonus damage...") needs
#    to end up collapsed right after this newly typed "duration" (i.e. this
#    was the last edited spot), so delete it from its old spot and re-add
#    it here.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(
    "as much stamina. Successful save causes them to loose half as much stamina or be exhausted",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$insPos2 = $rng2.Start + "as much stamina".Length
$ip2 = $d.Range($insPos2, $insPos2)
$ip2.InsertAfter(" as the duration")

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$goBackPos = $insPos2 + " as the duration".Length
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
